# Generate Report for Handback
# The handback for dea62695-3c19-4b4d-8586-a527fcf8cb49.md has now come in, so the
# report is regenerated: its status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the stale "handback file is not the latest"
# error clears, and the per-locale "Latest Handback DateTime" picks up the new
# handback timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-09-01 12:54:38"
$zhcn.Range("P3").Value = ""

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-09-01 12:54:45"
$dede.Range("P3").Value = ""

$zhcn.Columns.Item(16).ColumnWidth = 13.75
$dede.Columns.Item(16).ColumnWidth = 13.75
